$d = $word.ActiveDocument
$p = $d.Paragraphs(3)
$p4 = $d.Paragraphs(4)

# delete target, get insertion point
$target = $d.Range(30, 89)
$target.Delete()
$insPos = $p.Range.End - 1
$ins = $d.Range($insPos, $insPos)
$ins.InsertAfter("manter suas informacoes")
$newRange = $d.Range($insPos, $insPos + "manter suas informacoes".Length)

# try applying p4's formattedtext onto newRange
$donorRange = $d.Range($p4.Range.Start, $p4.Range.Start)
$newRange.FormattedText = $donorRange.FormattedText
Write-Output $p.Range.Text
